# Adds the vocabulary words and notes/quotes captured during the session.
# ENGLISH sheet gets 10 new word entries (rows 145-154).
# NOTES sheet gets 3 new quotes (rows 36-38).

$wb = $excel.ActiveWorkbook
$wsEnglish = $wb.Worksheets.Item("ENGLISH")
$wsNotes = $wb.Worksheets.Item("NOTES")

# --- ENGLISH: Word | Definition | Synonyms | Antonyms | Correct answer count | Created at | Tags
$englishRows = @(
    @("congested", "", "crowded", "", 0, "2021-11-20 22:32:31.286541", ""),
    @("swerve", "change or cause to change direction abruptly", "veer", "", 0, "2021-11-20 22:33:54.608594", ""),
    @("magnum opus", "a work of art, music, or literature that is regarded as the most important or best work that an artist, composer, or writer has produced", "", "", 0, "2021-11-20 22:36:19.119537", ""),
    @("anticipate", "", "expect;predict", "", 0, "2021-11-20 22:36:48.436207", ""),
    @("windfall", "a large amount of money that is won or received unexpectedly", "", "", 0, "2021-11-20 22:52:03.527958", ""),
    @("cushion", "", "pillow;protection", "", 0, "2021-11-20 22:53:43.81857", ""),
    @("diligently", "in a way that shows care in one's work or duties", "", "", 0, "2021-11-20 22:55:00.4924", ""),
    @("resilience", "", "flexibility", "", 0, "2021-11-20 22:56:23.511821", ""),
    @("parable", "", "allegory", "", 0, "2021-11-20 22:57:34.943717", ""),
    @("arbitrary", "", "random", "", 0, "2021-11-20 22:59:11.038374", "")
)

$startRow = 145
for ($i = 0; $i -lt $englishRows.Count; $i++) {
    $r = $startRow + $i
    $row = $englishRows[$i]
    $wsEnglish.Cells.Item($r, 1).Value = $row[0]
    $wsEnglish.Cells.Item($r, 2).Value = $row[1]
    $wsEnglish.Cells.Item($r, 3).Value = $row[2]
    $wsEnglish.Cells.Item($r, 4).Value = $row[3]
    $wsEnglish.Cells.Item($r, 5).Value = $row[4]
    $wsEnglish.Cells.Item($r, 6).Value = $row[5]
    $wsEnglish.Cells.Item($r, 7).Value = $row[6]
}

# --- NOTES: content | Tags
$noteRows = @(
    "Give me 6 hours to chop down a tree and I will spend the first 4 sharpening the axe",
    "The only thing we can expect (with any great certainty) is the unexpected",
    "To attain knowledge add things every day. To attain wisdom subtract things every day"
)

$startRow = 36
for ($i = 0; $i -lt $noteRows.Count; $i++) {
    $r = $startRow + $i
    $wsNotes.Cells.Item($r, 1).Value = $noteRows[$i]
    $wsNotes.Cells.Item($r, 2).Value = ""
}
